# Update "想去人数" (F) / "最低票价" (G) / Cover link (I) values on the
# "展览" and "全部类型" sheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Values shared by both "展览" (row 1 sheet) and "全部类型" sheet, except
# for row 31 column F which starts from a different value on each sheet
# but ends up at the same value (439) on both.
$commonUpdates = @(
    @{ Row = 2;  Col = "F"; Value = 19 }
    @{ Row = 4;  Col = "F"; Value = 264 }
    @{ Row = 6;  Col = "F"; Value = 543 }
    @{ Row = 7;  Col = "F"; Value = 52 }
    @{ Row = 8;  Col = "F"; Value = 2000 }
    @{ Row = 11; Col = "F"; Value = 4283 }
    @{ Row = 13; Col = "F"; Value = 281 }
    @{ Row = 15; Col = "F"; Value = 101 }
    @{ Row = 16; Col = "F"; Value = 21 }
    @{ Row = 19; Col = "F"; Value = 3037 }
    @{ Row = 20; Col = "F"; Value = 59 }
    @{ Row = 21; Col = "F"; Value = 439 }
    @{ Row = 22; Col = "G"; Value = 29.9 }
    @{ Row = 24; Col = "F"; Value = 70 }
    @{ Row = 25; Col = "F"; Value = 72 }
    @{ Row = 26; Col = "F"; Value = 8 }
    @{ Row = 26; Col = "G"; Value = 29.9 }
    @{ Row = 26; Col = "I"; Value = "//i0.hdslb.com/bfs/openplatform/202409/r7juTsXz1726729959726.jpeg" }
    @{ Row = 28; Col = "F"; Value = 48 }
    @{ Row = 28; Col = "G"; Value = 29.9 }
    @{ Row = 31; Col = "F"; Value = 439 }
    @{ Row = 32; Col = "F"; Value = 1685 }
    @{ Row = 33; Col = "F"; Value = 247 }
)

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($update in $commonUpdates) {
        $ws.Range("$($update.Col)$($update.Row)").Value = $update.Value
    }
}
